$d = $word.ActiveDocument

# Locate the final paragraph of the body (contains "Imagen 13" drawing + the _GoBack bookmark)
$lastP = $d.Paragraphs.Last
$insertionPoint = $lastP.Range.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

# The freshly-created (empty) paragraph is now the new last paragraph;
# replace its OOXML with the five target paragraphs in one shot.
$newP = $d.Paragraphs.Last
$xml = @'
<w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="113" w:after="113"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="273B47"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="273B47"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Formas de incluir media </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="273B47"/></w:rPr><w:t>queries</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="113" w:beforeAutospacing="0" w:after="113" w:afterAutospacing="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="273B47"/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="273B47"/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">En esta clase aprenderás a insertar un media </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="273B47"/><w:sz w:val="22"/></w:rPr><w:t>querie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="273B47"/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> en tu proyecto. Para ello, vas a trabajar sobre tu hoja de estilos, utilizando el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="273B47"/><w:sz w:val="22"/></w:rPr><w:t>tag</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="273B47"/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="273B47"/><w:sz w:val="22"/></w:rPr><w:t>style</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="273B47"/><w:sz w:val="22"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="113" w:beforeAutospacing="0" w:after="113" w:afterAutospacing="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="273B47"/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="273B47"/><w:sz w:val="22"/></w:rPr><w:t>El primer paso para lograr esto será realizar una nueva hoja de estilos en tu proyecto, ésta debe contar, en primer lugar, con la etiqueta link; harás uso de la aplicación de medidas para la pantalla, bordes y colores, entre otras características.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="113" w:beforeAutospacing="0" w:after="113" w:afterAutospacing="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="273B47"/><w:sz w:val="22"/><w:u w:val="single"/></w:rPr></w:pPr></w:p>
'@
$newP.Range.InsertXML($xml)

# Footer page-number field cached result: "3" -> "4" (one extra page of content).
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$ch = $ftr.Range.Characters.Item(1)
if ($ch.Text -eq "3") {
    $ch.Text = "4"
}

Write-Output "paragraphs=$($d.Paragraphs.Count)"
